$wb = $excel.ActiveWorkbook

# --- Overview sheet: Status column reused the same text for both locales ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B2").Value = "Handoff transform failed"
$wsOverview.Range("C2").Value = "Handoff transform failed"

# --- Per-locale detail sheets: Generate Report for handoff ---
$locales = @("zh-cn", "de-de")
foreach ($locale in $locales) {
    $ws = $wb.Worksheets.Item($locale)

    # Status -> Handoff transform failed
    $ws.Range("B2").Value = "Handoff transform failed"

    # Latest Handoff File: clear value + drop the hyperlink entirely
    $ws.Range("C2").Hyperlinks.Delete()
    $ws.Range("C2").Clear()

    # Latest Handoff Datetime: reset to the zero-value sentinel
    $ws.Range("D2").Value = "0001-01-01 00:00:00"

    # Handoff Reason: Include -> Ignored
    $ws.Range("H2").Value = "Ignored"
}
